{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph that holds the \"Charmander\" answer text - the last\n// paragraph of the \"How do we use Inheritance in our code?\" Q&A pair.\nconst marker = \"Charmander class inherits from a Pokemon class\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find the Charmander paragraph to anchor the new content on.\");\n}\n\n// New Q&A content about static variables/methods, inserted right after the\n// Inheritance section, following the doc's existing \"question / blank /\n// answer / blank\" paragraph pattern.\nconst newParagraphs = [\n  \"\",\n  \"What are static variables and methods?\",\n  \"\",\n  \"A Static Variable is a variable that is declared with the static keyword, meaning it retains its value across all instances of the class. A Static method is a method that is also declared with the static keyword, allowing it to be called without an object. It can only directly access static variables and other static methods within the class. Static variables and methods belong to the class itself rather than to any specific instance of that class. This means they are shared among all instances of the class. A static variable holds a value that is common to all objects of the class, while a static method can be called on the class itself without needing an object. \",\n  \"\",\n  \"Why do we use static variables and methods?\",\n  \"\",\n  \"Static variables and methods are useful when you want to share a single piece of data or behavior across all instances of a class, or when you want to perform operations that are independent of object state. They save memory because they are stored only once, and they can simplify access to shared resources or utility functions.\",\n  \"\",\n  \"How do we use static variables and methods in our code?\",\n  \"\",\n  \"To declare a static variable or method, you use the static keyword. For example totalrounds is a static variable shared across all instances of the Arena class.\"\n];\n\nlet insertAfter = anchor;\nfor (const text of newParagraphs) {\n  insertAfter = insertAfter.insertParagraph(text, \"After\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph holding the Inheritance section's last answer (the\n# \"Charmander\" paragraph) so the new static-variables Q&A content can be\n# inserted right after it, following the doc's existing\n# question / blank / answer / blank paragraph pattern.\n$marker = \"Charmander class inherits from a Pokemon class\"\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$marker*\") {\n        $anchor = $p\n        break\n    }\n}\n\nif ($anchor -eq $null) {\n    throw \"Could not find the Charmander paragraph to anchor the new content on.\"\n}\n\n$newParagraphs = @(\n    \"\",\n    \"What are static variables and methods?\",\n    \"\",\n    \"A Static Variable is a variable that is declared with the static keyword, meaning it retains its value across all instances of the class. A Static method is a method that is also declared with the static keyword, allowing it to be called without an object. It can only directly access static variables and other static methods within the class. Static variables and methods belong to the class itself rather than to any specific instance of that class. This means they are shared among all instances of the class. A static variable holds a value that is common to all objects of the class, while a static method can be called on the class itself without needing an object. \",\n    \"\",\n    \"Why do we use static variables and methods?\",\n    \"\",\n    \"Static variables and methods are useful when you want to share a single piece of data or behavior across all instances of a class, or when you want to perform operations that are independent of object state. They save memory because they are stored only once, and they can simplify access to shared resources or utility functions.\",\n    \"\",\n    \"How do we use static variables and methods in our code?\",\n    \"\",\n    \"To declare a static variable or method, you use the static keyword. For example totalrounds is a static variable shared across all instances of the Arena class.\"\n)\n\n$cur = $anchor\nforeach ($t in $newParagraphs) {\n    $cur.Range.InsertParagraphAfter()\n    $cur = $cur.Next()\n    if ($t -ne \"\") {\n        $cur.Range.InsertAfter($t)\n    }\n}\n"}
